$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "goibibo_Config" (sheet 1): update Execution / TestCase
# columns and add the Sl_no values for the two new rows.
# ---------------------------------------------------------------
$wsConfig = $wb.Worksheets.Item("goibibo_Config")

$wsConfig.Range("B2").Value = "No"
$wsConfig.Range("A3").Value = "2"
$wsConfig.Range("C3").Value = "searchHotelAndVerifyHotelDetail_TC02"
$wsConfig.Range("A4").Value = "3"
$wsConfig.Range("C4").Value = "searchHotelAndApplyFilter_TC03"

# Copy A2's formatting (bordered + quote-prefixed number style) down
# onto the two freshly populated cells in column A.
$wsConfig.Range("A2").Copy()
$wsConfig.Range("A3:A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# Sheet "TestDataSheet" (sheet 2): add two more test-data rows
# ---------------------------------------------------------------
$wsData = $wb.Worksheets.Item("TestDataSheet")

$wsData.Range("A3").Value = "searchHotelAndVerifyHotelDetail_TC02"
$wsData.Range("B3").Value = "edge"
$wsData.Range("C3").Value = "https://www.goibibo.com/"
$wsData.Range("D3").Value = "ooty"

$wsData.Range("A4").Value = "searchHotelAndApplyFilter_TC03"
$wsData.Range("B4").Value = "edge"
$wsData.Range("C4").Value = "https://www.goibibo.com/"
$wsData.Range("D4").Value = "ooty"

# Copy row 2's formatting down onto the two new rows.
$wsData.Range("A2:D2").Copy()
$wsData.Range("A3:D4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------
$wsConfig.Range("C4").Select()

$wsData.Activate()
$wsData.Range("F7").Select()

Write-Output "done"
